$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.233.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = "'3.562.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.08%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").Value = "'605.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.10%  '

$ws.Range("D6").Value = "'147.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.93%  '

$ws.Range("D7").Value = "'3.562.70"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("E8").Value = '  -0.19%  '

$ws.Range("E9").Value = '  -0.11%  '

$ws.Range("D10").Value = "'0.134"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '

$ws.Range("D11").Value = "'7.89"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.02%  '

$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("D13").Value = "'4.170.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.02%  '

$ws.Range("D14").Value = "'0.0000203"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("E15").Value = '  -3.15%  '

$ws.Range("D16").Value = "'3.565.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("E17").Value = '  +1.71%  '

$ws.Range("D18").Value = "'66.249.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.02%  '

$ws.Range("E19").Value = '  -4.04%  '

$ws.Range("D20").Value = "'6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.66%  '

$ws.Range("D21").Value = "'14.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.50%  '

$ws.Range("D22").Value = "'420.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.37%  '

$ws.Range("E23").Value = '  -0.98%  '

$ws.Range("D24").Value = "'77.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.24%  '

$ws.Range("D25").Value = "'3.704.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.10%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("E28").Value = '  +1.37%  '

$ws.Range("D29").Value = "'8.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.09%  '

$ws.Range("E30").Value = '  -0.56%  '

$ws.Range("D31").Value = "'1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.01%  '

$ws.Range("D32").Value = "'3.559.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.02%  '

$ws.Range("E33").Value = '  +3.92%  '

$ws.Range("D34").Value = "'24.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.72%  '

$ws.Range("E35").Value = '  -3.35%  '

$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").Value = "'7.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.61%  '

$ws.Range("E38").Value = '  -3.20%  '

$ws.Range("D39").Value = "'1.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.62%  '

$ws.Range("D40").Value = "'175.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.26%  '

$ws.Range("D41").Value = "'0.0841"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.11%  '

$ws.Range("D42").Value = "'5.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.95%  '

$ws.Range("E43").Value = '  -1.47%  '

$ws.Range("D44").Value = "'45.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.51%  '

$ws.Range("D45").Value = "'1.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.74%  '

$ws.Range("E46").Value = '  +0.07%  '

$ws.Range("D47").Value = "'2.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '

$ws.Range("D48").Value = "'23.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.33%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'24.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.43%  '

$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = "'7.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.29%  '

$ws.Range("B51").Value = 'ONDO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D51").Value = "'1.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.41%  '
